# Refresh the "cryptos" price/volume list (Price = column D, Volume(1h) = column E)
# for the crypto rows that moved since the last scrape. Values that would
# otherwise be auto-parsed by Excel as numbers (single "." decimal, e.g.
# "208.57") are written with a leading apostrophe so they stay plain text,
# exactly like the other Price cells (e.g. "27.042.87") that already resist
# numeric coercion because of their multi-dot thousands formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.042.87"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.563.23"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'208.57"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("D8").Value = "'22.08"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D12").Value = "1.785.95"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "1.568.12"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "27.043.90"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'61.90"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "'215.81"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'153.81"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "'0.106"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'3.20"
$ws.Range("E33").Value = "  +4.08%  "
$ws.Range("D34").Value = "1.425.01"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("E36").Value = "  +10.37%  "
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "'5.79"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'64.78"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "1.700.09"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").Value = "'86.74"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  +0.43%  "
